# Apply the updated "ReadMyAssessments - Back End - Team Lead" test case content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 : Step 1 ----
$ws.Range("C2").Value = 'Step 1: while logged out try to enter the "my assessments" page'
$ws.Range("D2").Value = "I am returned to the login screen"
$ws.Range("E2").Clear()
$ws.Range("F2").Value = "pass"

# ---- Row 3 : Step 2 ----
$ws.Range("C3").Value = "Step 2: Login as a user with the correct role"
$ws.Range("D3").Value = "I am redirected to the dashboard of the user"
$ws.Range("E3").Clear()
$ws.Range("F3").Value = "pass"

# ---- Row 4 : Step 3 ----
$ws.Range("C4").Value = ' Step 3: Click "My Assessments" '
$ws.Range("D4").Value = "I am showed a list of assessments all based around me"
$ws.Range("E4").Clear()
$ws.Range("F4").Value = "Pass"

# ---- Row 5 : Step 4 ----
$ws.Range("C5").Value = "Step 4: From the url view a list of someone elses assessments"
$ws.Range("D5").Value = "I am redirected to a list of my assessments(I should not be able to see a list of other people's assessments)"
$ws.Range("E5").Value = "directed me to someone else's assessment"
$ws.Range("F5").Value = "fail"

# ---- Row 6 : Step 5 (previously blank row) ----
$ws.Range("C6").Value = "Step 5: From the url change the myassessments id to nothing"
$ws.Range("D6").Value = "Redirected to a page saying access denied"
$ws.Range("E6").Value = "the application threw an exception "
$ws.Range("F6").Value = "fail"

# ---- Update the view: selected cell moves to F2, scrolled down a bit ----
$ws.Range("F2").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 2 | Out-Null
